$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCols = @(2, 6, 10, 14, 18)   # B, F, J, N, R
$newCols = @(1, 2, 3, 4, 5)      # A, B, C, D, E

for ($r = 4; $r -le 29; $r++) {
    $vals = @()
    foreach ($oc in $oldCols) {
        $vals += , ($ws.Cells.Item($r, $oc).Value())
    }
    for ($i = 0; $i -lt $newCols.Length; $i++) {
        $ws.Cells.Item($r, $newCols[$i]).Value = $vals[$i]
    }
}

# Clear stale old columns that are no longer part of the new layout (F, J, N, R)
$staleCols = @(6, 10, 14, 18)
foreach ($c in $staleCols) {
    $ws.Range($ws.Cells.Item(4, $c), $ws.Cells.Item(29, $c)).Clear()
}

$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H4").Select()
